$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.405737400054932
$ws.Range("B1").Value = 1.472282767295837
$ws.Range("C1").Value = 1.662174105644226
$ws.Range("D1").Value = 2.590127468109131
$ws.Range("E1").Value = -1
